# Update loading_percent results for the 380 kV case (Case_5_0)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> ordered list of [Column, NewValue] pairs
$updates = @{
    2 = @{ "B"=29.2202491786634; "C"=11.00672983195742; "D"=3.364191789476482; "E"=9.278120118369966; "F"=57.50620784017602; "J"=9.681245611534155; "L"=11.66435812606972; "M"=22.33059044669317; "N"=23.80895137667833 }
    3 = @{ "B"=28.98779264211935; "C"=10.67435376955665; "D"=3.304570384338744; "E"=9.263235039810452; "F"=57.42628700465458; "J"=9.690777016730125; "L"=11.68308512607782; "M"=22.30601294436343; "N"=23.86462593616423 }
    4 = @{ "B"=28.85223907291067; "C"=10.46883087925419; "D"=3.266873852802869; "E"=9.253889449250595; "F"=57.38978736585813; "J"=9.696921398371323; "L"=11.69613837047828; "M"=22.29553293496956; "N"=23.90079462990975 }
    5 = @{ "B"=28.79885519298272; "C"=10.38486423289609; "D"=3.251244784431305; "E"=9.250028918752069; "F"=57.37807621702649; "J"=9.699498972867117; "L"=11.70184883001646; "M"=22.29242433019906; "N"=23.91603279164929 }
    6 = @{ "B"=28.79010431185246; "C"=10.37091317956711; "D"=3.248633628645798; "E"=9.249384723398993; "F"=57.37632260536611; "J"=9.699931435075877; "L"=11.70282067743589; "M"=22.29197839929428; "N"=23.91859322781216 }
    7 = @{ "B"=28.85151154278653; "C"=10.46769914173208; "D"=3.266664147880525; "E"=9.253837596223777; "F"=57.38961661962541; "J"=9.696955861751109; "L"=11.69621379973292; "M"=22.29548630310006; "N"=23.90099811602152 }
    8 = @{ "B"=29.13864049072093; "C"=10.89251240622517; "D"=3.343864014671925; "E"=9.273030333322662; "F"=57.47604333462374; "J"=9.684471572823645; "L"=11.67049262073112; "M"=22.32116101448032; "N"=23.82773621999416 }
    9 = @{ "B"=29.75623174826058; "C"=11.70804933666505; "D"=3.486352041424167; "E"=9.309046628936315; "F"=57.74512274820415; "J"=9.662295882805154; "L"=11.63238140752484; "M"=22.40794414869565; "N"=23.69980657159574 }
    10 = @{ "B"=30.23979210309301; "C"=12.2890918727344; "D"=3.585286108132792; "E"=9.334541539041604; "F"=58.00318329665982; "J"=9.647393097499915; "L"=11.61188460691332; "M"=22.49367241881273; "N"=23.6153988050694 }
    11 = @{ "B"=30.46547013018463; "C"=12.54811432529462; "D"=3.62898238061793; "E"=9.345934723957434; "F"=58.13356501851217; "J"=9.640911753624172; "L"=11.60418633049808; "M"=22.53737493345046; "N"=23.57907828992031 }
    12 = @{ "B"=30.55168059710588; "C"=12.64533743405313; "D"=3.64533586842529; "E"=9.350219944382951; "F"=58.18479000950693; "J"=9.638500029779866; "L"=11.60150466037358; "M"=22.554593323555; "N"=23.56562331700719 }
    13 = @{ "B"=30.53308133786409; "C"=12.62443859125059; "D"=3.641822528191192; "E"=9.349298336576318; "F"=58.17367568376153; "J"=9.639017545800272; "L"=11.60207182604459; "M"=22.55085539640962; "N"=23.56850779343933 }
    14 = @{ "B"=30.47254806408145; "C"=12.55613077130127; "D"=3.6303316945585; "E"=9.346287853819515; "F"=58.13774230465785; "J"=9.6407124868048; "L"=11.6039610299865; "M"=22.538778133462; "N"=23.57796535149251 }
    15 = @{ "B"=30.43556542960567; "C"=12.51417502408223; "D"=3.623267903843314; "E"=9.344440060920576; "F"=58.11597279022897; "J"=9.641756230657652; "L"=11.60514862131356; "M"=22.53146738133476; "N"=23.58379729870336 }
    16 = @{ "B"=30.22515241031901; "C"=12.27204866277476; "D"=3.582403666449474; "E"=9.333792875516099; "F"=57.99492264131073; "J"=9.647822645323457; "L"=11.61242039898456; "M"=22.49091041469273; "N"=23.61781424262064 }
    17 = @{ "B"=30.0974829828116; "C"=12.12207998800809; "D"=3.556995720617052; "E"=9.3272089792359; "F"=57.92397876272241; "J"=9.651620355230087; "L"=11.61729761487996; "M"=22.4672299793775; "N"=23.63921460068998 }
    18 = @{ "B"=30.02459202233211; "C"=12.03532883014758; "D"=3.542258783062229; "E"=9.323402823371588; "F"=57.88439736140388; "J"=9.653832761289461; "L"=11.62025591009094; "M"=22.45405289460662; "N"=23.65171900709742 }
    19 = @{ "B"=30.00000741902501; "C"=12.00587493721649; "D"=3.537248150674446; "E"=9.322110805521421; "F"=57.87120637239678; "J"=9.654586671052662; "L"=11.62128383271306; "M"=22.44966770212439; "N"=23.65598635520192 }
    20 = @{ "B"=30.1110181047247; "C"=12.13809621457026; "D"=3.559713203364967; "E"=9.327911838821512; "F"=57.93140432200812; "J"=9.651213179724701; "L"=11.61676258886759; "M"=22.4697049751166; "N"=23.63691626297423 }
    21 = @{ "B"=30.49030832592218; "C"=12.57621861565169; "D"=3.633712116859988; "E"=9.347172893922671; "F"=58.14824667478175; "J"=9.640213486689143; "L"=11.60339979036239; "M"=22.54230741830704; "N"=23.57517932626459 }
    22 = @{ "B"=30.7425416108897; "C"=12.85748626815443; "D"=3.6809456285757; "E"=9.359591438815322; "F"=58.30075320187805; "J"=9.633272868601576; "L"=11.59602728615274; "M"=22.59365399215803; "N"=23.53657249984476 }
    23 = @{ "B"=30.60754548081956; "C"=12.70786331414567; "D"=3.655841108120513; "E"=9.352978857944871; "F"=58.2183762834043; "J"=9.636954562669475; "L"=11.59983771714636; "M"=22.5658954120829; "N"=23.55701825493139 }
    24 = @{ "B"=30.10489728890272; "C"=12.1308569283258; "D"=3.558485032611156; "E"=9.327594141391311; "F"=57.92804347164397; "J"=9.651397173262382; "L"=11.61700399343726; "M"=22.46858466749114; "N"=23.63795471522349 }
    25 = @{ "B"=29.58366235797698; "C"=11.49010924548687; "D"=3.448792992601379; "E"=9.299474512422909; "F"=57.66169323634106; "J"=9.668049792046599; "L"=11.64137281694567; "M"=22.3805883422199; "N"=23.73273149448165 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

Write-Output "Applied 380 kV case updates to loading_percent sheet"